$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. Insert a new row at
# position 14, which pushes the existing rows 14-18 down to rows 15-19,
# preserving their data and formatting (date format, etc.).
$ws.Rows("14").Insert()

# Fill in the newly inserted row 14 with the new record's values.
$ws.Cells.Item(14, 1).Value2 = 10
$ws.Cells.Item(14, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(14, 3).Value2 = "La Araucanía"
$ws.Cells.Item(14, 4).Value2 = 44902
$ws.Cells.Item(14, 4).NumberFormat = $ws.Cells.Item(15, 4).NumberFormat
$ws.Cells.Item(14, 5).Value2 = 9
$ws.Cells.Item(14, 6).Value2 = "Fruta"
$ws.Cells.Item(14, 7).Value2 = 100104
$ws.Cells.Item(14, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(14, 9).Value2 = 100104004
$ws.Cells.Item(14, 10).Value2 = "Níspero"
$ws.Cells.Item(14, 11).Value2 = "Californiana(o)"
$ws.Cells.Item(14, 12).Value2 = "Primera"
$ws.Cells.Item(14, 13).Value2 = 90
$ws.Cells.Item(14, 14).Value2 = 25000
$ws.Cells.Item(14, 15).Value2 = 25000
$ws.Cells.Item(14, 16).Value2 = 25000
$ws.Cells.Item(14, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(14, 18).Value2 = "Provincia de Quillota"
$ws.Cells.Item(14, 19).Value2 = 2500
$ws.Cells.Item(14, 20).Value2 = 10
